# Weekly fruit/vegetable update: a new week's worth of prices (Primera +
# Segunda quality) is inserted for "Zanahoria" at Terminal Hortofrutícola
# Agro Chillán, pushing the existing rows 435:486 down to 437:488.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 435; this shifts
# the old rows 435:486 down to 437:488 and grows the used range to R488.
$ws.Rows.Item(435).Insert()
$ws.Rows.Item(435).Insert()

# --- New row 435 (Primera quality, week of 2023-07-17) -------------------
$ws.Cells.Item(435, 1).Value = 7
$ws.Cells.Item(435, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(435, 3).Value = "Ñuble"
$ws.Cells.Item(435, 4).Value = 45124
$ws.Cells.Item(435, 5).Value = 16
$ws.Cells.Item(435, 6).Value = 100114013
$ws.Cells.Item(435, 7).Value = "Zanahoria"
$ws.Cells.Item(435, 8).Value = "Sin especificar"
$ws.Cells.Item(435, 9).Value = "Primera"
$ws.Cells.Item(435, 10).Value = 150
$ws.Cells.Item(435, 11).Value = 7000
$ws.Cells.Item(435, 12).Value = 7000
$ws.Cells.Item(435, 13).Value = 7000
$ws.Cells.Item(435, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(435, 15).Value = "Región de Ñuble"
$ws.Cells.Item(435, 16).Value = 350
$ws.Cells.Item(435, 17).Value = 20
$ws.Cells.Item(435, 18).Value = "Hortaliza"

# --- New row 436 (Segunda quality, same week) -----------------------------
$ws.Cells.Item(436, 1).Value = 7
$ws.Cells.Item(436, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(436, 3).Value = "Ñuble"
$ws.Cells.Item(436, 4).Value = 45124
$ws.Cells.Item(436, 5).Value = 16
$ws.Cells.Item(436, 6).Value = 100114013
$ws.Cells.Item(436, 7).Value = "Zanahoria"
$ws.Cells.Item(436, 8).Value = "Sin especificar"
$ws.Cells.Item(436, 9).Value = "Segunda"
$ws.Cells.Item(436, 10).Value = 100
$ws.Cells.Item(436, 11).Value = 6000
$ws.Cells.Item(436, 12).Value = 6000
$ws.Cells.Item(436, 13).Value = 6000
$ws.Cells.Item(436, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(436, 15).Value = "Región de Ñuble"
$ws.Cells.Item(436, 16).Value = 300
$ws.Cells.Item(436, 17).Value = 20
$ws.Cells.Item(436, 18).Value = "Hortaliza"
